$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38
$ws.Cells.Item($row, 1).Value = "Daniele Ruzzenenti "
$ws.Cells.Item($row, 2).Value = "ELIA BATTISTI | U.S. Guarna"
$ws.Cells.Item($row, 3).Value = "Carlo Stedile | MAI UNA GIOIA"
$ws.Cells.Item($row, 4).Value = "Alessandro  Maffei | FC Savignano"
$ws.Cells.Item($row, 5).Value = "Riccardo Zeni | Demobusters"
$ws.Cells.Item($row, 6).Value = "Jacopo Zecchinelli | Vigili del Fusto"
